# Update workbook "广州-漫展信息.xlsx" to the scraped data as of commit 456a3b4.
#
# The underlying dataset removed a duplicated "广州·Look Look动漫嘉年华" entry
# (it had been scraped twice, back to back, in both the "展览" sheet and the
# combined "全部类型" sheet) and several "想去人数" (interest-count) /
# "最低票价" (min-price) values were refreshed from a later scrape.
#
# Concretely:
#   - Sheet "展览"   (sheet index 1): delete the duplicate row 4, then
#       refresh a handful of F/G values in the rows that shifted up.
#   - Sheet "演出"   (sheet index 2): no rows added/removed, just F7 refreshed.
#   - Sheet "本地生活"(sheet index 3): no rows added/removed, just F2 refreshed.
#   - Sheet "全部类型"(sheet index 4): same duplicate-row removal (at row 5
#       here, since this sheet carries two extra leading rows versus
#       "展览"), plus the same F/G refreshes (offset by one row) and the
#       two standalone F refreshes that also appear on "演出"/"本地生活".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Remove the duplicated "Look Look" row (old row 4); everything below
# shifts up by one, and the sheet's used range shrinks from A1:I25 to
# A1:I24 automatically.
$ws1.Rows.Item(4).Delete()

# Refresh counts on the rows that shifted up (new row numbers).
$ws1.Range("F3").Value = 1814
$ws1.Range("F4").Value = 399
$ws1.Range("F5").Value = 1482
$ws1.Range("F8").Value = 734
$ws1.Range("F9").Value = 13155
$ws1.Range("G9").Value = 78
$ws1.Range("F10").Value = 13028
$ws1.Range("F15").Value = 65
$ws1.Range("F16").Value = 631
$ws1.Range("F17").Value = 2058
$ws1.Range("F19").Value = 27
$ws1.Range("F22").Value = 188
$ws1.Range("F23").Value = 272
$ws1.Range("F24").Value = 735

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performances) - no structural change, one refreshed count
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F7").Value = 102

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life) - no structural change, one refreshed count
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 186

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types, combined / chronological view)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# This sheet has its own copy of the duplicated "Look Look" row, one row
# further down (row 5, because rows 2-3 here hold entries that live on
# other sheets). Remove it the same way.
$ws4.Rows.Item(5).Delete()

# Refresh counts (rows shifted up by one relative to before).
$ws4.Range("F2").Value = 186
$ws4.Range("F4").Value = 1814
$ws4.Range("F5").Value = 399
$ws4.Range("F6").Value = 1482
$ws4.Range("F10").Value = 734
$ws4.Range("F11").Value = 13155
$ws4.Range("G11").Value = 78
$ws4.Range("F12").Value = 13028
$ws4.Range("F17").Value = 65
$ws4.Range("F18").Value = 631
$ws4.Range("F21").Value = 2058
$ws4.Range("F23").Value = 27
$ws4.Range("F28").Value = 188
$ws4.Range("F29").Value = 272
$ws4.Range("F30").Value = 735
$ws4.Range("F31").Value = 102
